$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("June 10, 2022", $true, $false, $false, $false, $false, $true, 1, $false, "June 11, 2022", 2)

$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Replacement.ClearFormatting()
$find2.Execute("August 09, 2022", $true, $false, $false, $false, $false, $true, 1, $false, "August 10, 2022", 2)
